# Edit: replace the textual NPS formula answer in B2 with the actual
# numeric result (0.2), displayed as a percentage ("20%") via a custom
# number format, matching the "((Promotores - Detratores) / Total de
# Avaliações) * 100 = 20%" text that previously lived in that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("B2").NumberFormat = "#,##0%"
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 11
